$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 11146.111
$ws.Range("I34").Value = 759.2857
$ws.Range("J34").Value = 47500
$ws.Range("K34").Value = 759.2857
$ws.Range("L34").Value = 47500
$ws.Range("M34").Value = -556.2857
$ws.Range("N34").Value = -47906

# Row 36
$ws.Range("H36").Value = 11146.111
$ws.Range("I36").Value = 759.2857
$ws.Range("J36").Value = 47500
$ws.Range("K36").Value = 759.2857
$ws.Range("L36").Value = 47500
$ws.Range("M36").Value = -44.28570000000002
$ws.Range("N36").Value = -48930

# Row 76
$ws.Range("H76").Value = 10591.412
$ws.Range("I76").Value = 12815.692
$ws.Range("K76").Value = 12815.692
$ws.Range("M76").Value = -12500.692

# Row 79
$ws.Range("H79").Value = 10591.412
$ws.Range("I79").Value = 12815.692
$ws.Range("K79").Value = 12815.692
$ws.Range("M79").Value = -11723.692

# Row 112
$ws.Range("H112").Value = 4482.5
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 4980
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 14940
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -17156

# Row 115
$ws.Range("H115").Value = 663.625
$ws.Range("I115").Value = 329.85715
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 989.5714499999999
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = 577.4285500000001
$ws.Range("N115").Value = -12134

# Row 118
$ws.Range("H118").Value = 446164.22
$ws.Range("I118").Value = 2001000
$ws.Range("J118").Value = 1925.4286
$ws.Range("K118").Value = 6003000
$ws.Range("L118").Value = 5776.2858
$ws.Range("M118").Value = -6001343
$ws.Range("N118").Value = -9090.2858

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 7302
$ws.Range("I25").Value = 7302
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 7302
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -6900
$ws.Range("N25").ClearContents()

# Row 38
$ws.Range("H38").Value = 3020
$ws.Range("I38").Value = 3019
$ws.Range("J38").Value = 3021
$ws.Range("K38").Value = 3019
$ws.Range("L38").Value = 3021
$ws.Range("M38").Value = -2552
$ws.Range("N38").Value = -3955

# Row 74
$ws.Range("H74").Value = 1556.5238
$ws.Range("I74").Value = 1013.1667
$ws.Range("J74").Value = 4816.6665
$ws.Range("K74").Value = 1013.1667
$ws.Range("L74").Value = 4816.6665
$ws.Range("M74").Value = -139.1667
$ws.Range("N74").Value = -6564.6665

# Row 77
$ws.Range("H77").Value = 1556.5238
$ws.Range("I77").Value = 1013.1667
$ws.Range("J77").Value = 4816.6665
$ws.Range("K77").Value = 5065.8335
$ws.Range("L77").Value = 24083.3325
$ws.Range("M77").Value = -697.8334999999997
$ws.Range("N77").Value = -32819.3325

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 675.9231
$ws.Range("I94").Value = 646.5625
$ws.Range("J94").Value = 722.9
$ws.Range("K94").Value = 646.5625
$ws.Range("L94").Value = 722.9
$ws.Range("M94").Value = -195.5625
$ws.Range("N94").Value = -1624.9

# Row 118
$ws.Range("H118").Value = 14900
$ws.Range("J118").Value = 14900
$ws.Range("L118").Value = 14900
$ws.Range("N118").Value = -18214

# Row 122
$ws.Range("H122").Value = 33405.555
$ws.Range("J122").Value = 33405.555
$ws.Range("L122").Value = 33405.555
$ws.Range("N122").Value = -43205.555

# Row 134
$ws.Range("H134").Value = 1629.6666
$ws.Range("I134").Value = 1644.5
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 4933.5
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -2398.5
$ws.Range("N134").Value = -9870

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3830.2307
$ws.Range("I16").Value = 4143.636
$ws.Range("J16").Value = 2106.5
$ws.Range("K16").Value = 4143.636
$ws.Range("L16").Value = 2106.5
$ws.Range("M16").Value = -3856.636
$ws.Range("N16").Value = -2680.5

# Row 22
$ws.Range("H22").Value = 231
$ws.Range("I22").Value = 235.71428
$ws.Range("J22").Value = 222.75
$ws.Range("K22").Value = 235.71428
$ws.Range("L22").Value = 222.75
$ws.Range("M22").Value = 114.28572
$ws.Range("N22").Value = -922.75

# Row 109
$ws.Range("H109").Value = 10857.143
$ws.Range("J109").Value = 10857.143
$ws.Range("L109").Value = 10857.143
$ws.Range("N109").Value = -12937.143

# Row 113
$ws.Range("H113").Value = 3830.2307
$ws.Range("I113").Value = 4143.636
$ws.Range("J113").Value = 2106.5
$ws.Range("K113").Value = 4143.636
$ws.Range("L113").Value = 2106.5
$ws.Range("M113").Value = -1973.636
$ws.Range("N113").Value = -6446.5

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

# Row 121
$ws.Range("H121").Value = 54980
$ws.Range("J121").Value = 54980
$ws.Range("L121").Value = 54980
$ws.Range("N121").Value = -57600

$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 290.44446
$ws.Range("I98").Value = 256
$ws.Range("J98").Value = 333.5
$ws.Range("K98").Value = 768
$ws.Range("L98").Value = 1000.5
$ws.Range("M98").Value = 730
$ws.Range("N98").Value = -3996.5

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2609.625
$ws.Range("I97").Value = 2551.111
$ws.Range("J97").Value = 2785.1667
$ws.Range("K97").Value = 2551.111
$ws.Range("L97").Value = 2785.1667
$ws.Range("M97").Value = -2055.111
$ws.Range("N97").Value = -3777.1667

# Row 132
$ws.Range("H132").Value = 4999.25
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Range("H18").Value = 15000
$ws.Range("J18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15344

# Row 20
$ws.Range("H20").Value = 37502.5
$ws.Range("I20").Value = 37005
$ws.Range("J20").Value = 38000
$ws.Range("K20").Value = 37005
$ws.Range("L20").Value = 38000
$ws.Range("M20").Value = -36779
$ws.Range("N20").Value = -38452

# Row 48
$ws.Range("H48").Value = 28046
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

# Row 68
$ws.Range("H68").Value = 1460.4814
$ws.Range("I68").Value = 1575
$ws.Range("J68").Value = 1393.1177
$ws.Range("K68").Value = 1575
$ws.Range("L68").Value = 1393.1177
$ws.Range("M68").Value = -826
$ws.Range("N68").Value = -2891.1177

# Row 71
$ws.Range("H71").Value = 1460.4814
$ws.Range("I71").Value = 1575
$ws.Range("J71").Value = 1393.1177
$ws.Range("K71").Value = 7875
$ws.Range("L71").Value = 6965.5885
$ws.Range("M71").Value = -4131
$ws.Range("N71").Value = -14453.5885

# Row 132
$ws.Range("H132").Value = 3514.2942
$ws.Range("I132").Value = 3537.4092
$ws.Range("J132").Value = 3471.9167
$ws.Range("K132").Value = 10612.2276
$ws.Range("L132").Value = 10415.7501
$ws.Range("M132").Value = -8082.2276
$ws.Range("N132").Value = -15475.7501

$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 17222.334
$ws.Range("J109").Value = 17222.334
$ws.Range("L109").Value = 17222.334
$ws.Range("N109").Value = -19996.334
